$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7263558506965637
$ws.Range("B1").Value = 1.027143597602844
$ws.Range("C1").Value = 0.9339823126792908
$ws.Range("D1").Value = 3.249042749404907
$ws.Range("E1").Value = 1.614280343055725
